$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Plasma_Gen" updates
# ---------------------------------------------------------------------------
$wsGen = $wb.Worksheets.Item("Plasma_Gen")

# Row 5 / H5 - new answer text
$wsGen.Range("H5").Value = "Straight type or Right angle type?"

# Row 9 - taller row + new answer (contact info block)
$wsGen.Rows.Item(9).RowHeight = 82.5
$wsGen.Range("H9").Value = "GSP Korea`n유재희`nMobile: +82 10 8648 2090`nFax     : +82 31 427 8523`nE-mail: jhyoo@gspkorea.co.kr"

# Row 10 - new answer (proposal text)
$wsGen.Range("H10").Value = "제안 : USB CON 하나로 RS-232 통신과 External 전원 공급 지원"

# Row 17 - brand new Q&A entry
$wsGen.Rows.Item(17).RowHeight = 99

$wsGen.Range("C17").Value = 43132

$wsGen.Range("D17").Value = "H/W"
$wsGen.Range("D17").HorizontalAlignment = -4108
$wsGen.Range("D17").VerticalAlignment = -4108

$wsGen.Range("E17").Value = "OPEN"
$wsGen.Range("E17").HorizontalAlignment = -4108
$wsGen.Range("E17").VerticalAlignment = -4108

$wsGen.Range("F17").Value = "Certification을 받아야 하는 인증 목록 확인 필요`nBattery를 사용하는 경우 추가 인증이 있을 수 있음"
$wsGen.Range("F17").WrapText = $true

$wsGen.Range("H17").Value = "보통 국가마다 제품 인증이 필요하며, 아래는 전자제품의 경우 받는 인증들임.`n장비의 경우는 제가 잘 모릅니다.  `n  국내 : KC 인증`n  유럽 : EC 인증`n  미국 : FCC 인증 등등.."

# Column F width shrinks a bit (no longer "best fit")
$wsGen.Columns.Item(6).ColumnWidth = 69.43

# Selection moves to G21
$wsGen.Range("G21").Select()

# ---------------------------------------------------------------------------
# Sheet "Plasma_LF" updates
# ---------------------------------------------------------------------------
$wsLF = $wb.Worksheets.Item("Plasma_LF")

# Column G narrower
$wsLF.Columns.Item(7).ColumnWidth = 37.85

# Print scale 66% -> 68% (keep "fit to page" behaviour / fitToHeight = 0)
$wsLF.PageSetup.Zoom = 68
$wsLF.PageSetup.FitToPagesTall = $false

# Selection moves to F12 - Plasma_LF stays the active/visible tab
$wsLF.Range("F12").Select()
